## Generate Report for handoff
## Adds a new "handoff failed" row (55f27676-df51-4893-a5c3-80b71f5fab21.md)
## ahead of the existing ".localization-config" row on every sheet, and
## rewrites the 2e3a0a07... handoff (file + xlf + timestamps) to the new
## 190c0f8a... handoff on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newMd        = "190c0f8a-2aa5-4f77-8caf-e44bc19e3c7f.md"
$failedMd     = "55f27676-df51-4893-a5c3-80b71f5fab21.md"
$zhXlf        = "190c0f8a-2aa5-4f77-8caf-e44bc19e3c7f.06f7a964070a18f7a86a6c1306093ad89331d3af.zh-cn.xlf"
$deXlf        = "190c0f8a-2aa5-4f77-8caf-e44bc19e3c7f.06f7a964070a18f7a86a6c1306093ad89331d3af.de-de.xlf"
$zhHandoffDt  = "2016-01-08 15:45:55"
$deHandoffDt  = "2016-01-08 15:46:11"
$epoch        = "0001-01-01 00:00:00"

$mdBase    = "https://github.com/OpenLocalizationTest/oltest/blob/1a076aaffc3f9e05fa0a29fe74395b4e8c6547da/e2e/"
$cfgUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/1a076aaffc3f9e05fa0a29fe74395b4e8c6547da/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b8399e67d248ce7105952aa421fe36d64f202167/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5535b955b458f76cd1422cc5a0d240c9cbc9a6f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/"

## ---------------------------------------------------------------
## Sheet "Overview": File Name | zh-cn | de-de
## ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A1:C3").Hyperlinks.Delete()

$ov.Range("A2").Value = $newMd
$ov.Range("B2").Value = "Not yet handed off"
$ov.Range("C2").Value = "Not yet handed off"

$ov.Range("A3").Value = $failedMd
$ov.Range("B3").Value = "Handoff failed"
$ov.Range("C3").Value = "Handoff failed"

$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not localized"
$ov.Range("C4").Value = "Not localized"

$ov.Hyperlinks.Add($ov.Cells.Item(2,1), ($mdBase + $newMd), [System.Type]::Missing, [System.Type]::Missing, $newMd)
$ov.Hyperlinks.Add($ov.Cells.Item(3,1), ($mdBase + $failedMd), [System.Type]::Missing, [System.Type]::Missing, $failedMd)
$ov.Hyperlinks.Add($ov.Cells.Item(4,1), $cfgUrl, [System.Type]::Missing, [System.Type]::Missing, ".localization-config")

## ---------------------------------------------------------------
## Sheet "zh-cn"
## ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A1:I3").Hyperlinks.Delete()

$zh.Range("A2").Value = $newMd
$zh.Range("B2").Value = "Not yet handed off"
$zh.Range("C2").Value = $zhXlf
$zh.Range("D2").Value = $zhHandoffDt
$zh.Range("G2").Value = $epoch
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = $failedMd
$zh.Range("B3").Value = "Handoff failed"
$zh.Range("D3").Value = $epoch
$zh.Range("G3").Value = $epoch
$zh.Range("H3").Value = "Ignored"

$zh.Range("A4").Value = ".localization-config"
$zh.Range("B4").Value = "Not localized"
$zh.Range("D4").Value = $epoch
$zh.Range("G4").Value = $epoch
$zh.Range("H4").Value = "Ignored"

$zh.Hyperlinks.Add($zh.Cells.Item(2,1), ($mdBase + $newMd), [System.Type]::Missing, [System.Type]::Missing, $newMd)
$zh.Hyperlinks.Add($zh.Cells.Item(2,3), ($zhXlfBase + $zhXlf), [System.Type]::Missing, [System.Type]::Missing, $zhXlf)
$zh.Hyperlinks.Add($zh.Cells.Item(3,1), ($mdBase + $failedMd), [System.Type]::Missing, [System.Type]::Missing, $failedMd)
$zh.Hyperlinks.Add($zh.Cells.Item(4,1), $cfgUrl, [System.Type]::Missing, [System.Type]::Missing, ".localization-config")

## ---------------------------------------------------------------
## Sheet "de-de"
## ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A1:I3").Hyperlinks.Delete()

$de.Range("A2").Value = $newMd
$de.Range("B2").Value = "Not yet handed off"
$de.Range("C2").Value = $deXlf
$de.Range("D2").Value = $deHandoffDt
$de.Range("G2").Value = $epoch
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = $failedMd
$de.Range("B3").Value = "Handoff failed"
$de.Range("D3").Value = $epoch
$de.Range("G3").Value = $epoch
$de.Range("H3").Value = "Ignored"

$de.Range("A4").Value = ".localization-config"
$de.Range("B4").Value = "Not localized"
$de.Range("D4").Value = $epoch
$de.Range("G4").Value = $epoch
$de.Range("H4").Value = "Ignored"

$de.Hyperlinks.Add($de.Cells.Item(2,1), ($mdBase + $newMd), [System.Type]::Missing, [System.Type]::Missing, $newMd)
$de.Hyperlinks.Add($de.Cells.Item(2,3), ($deXlfBase + $deXlf), [System.Type]::Missing, [System.Type]::Missing, $deXlf)
$de.Hyperlinks.Add($de.Cells.Item(3,1), ($mdBase + $failedMd), [System.Type]::Missing, [System.Type]::Missing, $failedMd)
$de.Hyperlinks.Add($de.Cells.Item(4,1), $cfgUrl, [System.Type]::Missing, [System.Type]::Missing, ".localization-config")
